# Recomputed Step3 "First Rise Point" detection results after adding the
# configurable zero_before_threshold parameter (dims before the noise
# threshold / first rise point can now be forced to 0 before the pulse
# scan runs). This shifts First_Noticeable_Increase_Index (col C) and its
# matching First_Noticeable_Increase_Cumulative_Value (col E) for each
# signal segment, which in turn shifts the derived Pulse_Width (col G =
# Point_Exceeds_Index[D] - First_Noticeable_Increase_Index[C]) on every
# Step3_DataPts_* sheet (one per Intensity_Threshold: 0.5 / 0.7 / 0.8 / 0.9).

$wb = $excel.ActiveWorkbook

function Set-Step3Row {
    param($ws, $row, $firstIncreaseIndex, $firstIncreaseCumValue, $pulseWidth)
    $ws.Cells.Item($row, 3).Value = $firstIncreaseIndex   # C: First_Noticeable_Increase_Index
    $ws.Cells.Item($row, 5).Value = $firstIncreaseCumValue # E: First_Noticeable_Increase_Cumulative_Value
    $ws.Cells.Item($row, 7).Value = $pulseWidth            # G: Pulse_Width
}

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
Set-Step3Row $ws 2 87 0.004319649857380383 16
Set-Step3Row $ws 3 88 0.007759161124024368 17
Set-Step3Row $ws 4 87 0.006349619342386191 19
Set-Step3Row $ws 5 88 0.007763863460492989 19
Set-Step3Row $ws 6 88 0.009960658459728077 19

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
Set-Step3Row $ws 2 87 0.004319649857380383 29
Set-Step3Row $ws 3 88 0.007759161124024368 28
Set-Step3Row $ws 4 87 0.006349619342386191 31
Set-Step3Row $ws 5 88 0.007763863460492989 31
Set-Step3Row $ws 6 88 0.009960658459728077 29

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
Set-Step3Row $ws 2 87 0.004319649857380383 64
Set-Step3Row $ws 3 88 0.007759161124024368 64
Set-Step3Row $ws 4 87 0.006349619342386191 66
Set-Step3Row $ws 5 88 0.007763863460492989 65
Set-Step3Row $ws 6 88 0.009960658459728077 65

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
Set-Step3Row $ws 2 87 0.004319649857380383 75
Set-Step3Row $ws 3 88 0.007759161124024368 77
Set-Step3Row $ws 4 87 0.006349619342386191 81
Set-Step3Row $ws 5 88 0.007763863460492989 78
Set-Step3Row $ws 6 88 0.009960658459728077 80
